# Update gh-pages to output generated at 456a3b4
# Refresh the "想去人数" (F column) and one "最低票价" (G4) value on the
# "展览" and "全部类型" sheets.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 7344
$ws1.Range("F3").Value = 468
$ws1.Range("G4").Value = 40
$ws1.Range("F6").Value = 211
$ws1.Range("F7").Value = 146
$ws1.Range("F9").Value = 16
$ws1.Range("F10").Value = 70
$ws1.Range("F11").Value = 238
$ws1.Range("F14").Value = 35
$ws1.Range("F15").Value = 1892
$ws1.Range("F17").Value = 58
$ws1.Range("F18").Value = 3888
$ws1.Range("F25").Value = 2589
$ws1.Range("F27").Value = 345
$ws1.Range("F32").Value = 36
$ws1.Range("F36").Value = 171
$ws1.Range("F37").Value = 77
$ws1.Range("F38").Value = 1530
$ws1.Range("F39").Value = 188

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 7344
$ws4.Range("F3").Value = 468
$ws4.Range("G4").Value = 40
$ws4.Range("F7").Value = 211
$ws4.Range("F8").Value = 146
$ws4.Range("F10").Value = 16
$ws4.Range("F11").Value = 70
$ws4.Range("F12").Value = 238
$ws4.Range("F15").Value = 35
$ws4.Range("F16").Value = 1892
$ws4.Range("F18").Value = 58
$ws4.Range("F19").Value = 3888
$ws4.Range("F26").Value = 2589
$ws4.Range("F28").Value = 345
$ws4.Range("F33").Value = 36
$ws4.Range("F37").Value = 171
$ws4.Range("F38").Value = 77
$ws4.Range("F39").Value = 1530
$ws4.Range("F40").Value = 188
